$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Datos")

# --- Grow the sample-data table from 3 rows (2-4) to 5 rows (2-6) -----
# Copy the formatting of the last existing data row (row 4) down into the
# two new rows so they inherit the same cell styles as the existing ones.
$ws1.Range("A4:J4").Copy() | Out-Null
$ws1.Range("A5:J5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws1.Range("A6:J6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Remove the now-obsolete sample values in columns C:J for the original
# data rows (2-4); idCaso (A) / orientacion (B) are kept.
$ws1.Range("C2:E4").ClearContents()
$ws1.Range("F2:I4").ClearContents()
$ws1.Range("J2:J4").ClearContents()

# Populate the two additional test cases.
$ws1.Range("A5").Value = "'4"
$ws1.Range("B5").Value = "'Alterno"
$ws1.Range("A6").Value = "'5"
$ws1.Range("B6").Value = "'Alterno"

# Update selection to match the latest edit location.
$ws1.Activate() | Out-Null
$ws1.Range("J24").Select() | Out-Null
